# The presentation ships with two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (used by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" / "Red Violet" (used by the Slide
#                            Master, i.e. the presentation's design/theme)
#
# The authored edit swaps the two themes' contents, so that the deck's
# design (Slide Master) now uses the plain "Office" color scheme, while
# the (font/format scheme of both themes is identical already - only the
# color scheme differs between them). We reproduce the visible effect of
# that swap by rewriting the Slide Master / presentation theme's color
# scheme from the "Red Violet" values to the "Office" values.

function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target ("Office") color scheme values - these are the colors theme1.xml
# (the "Office Theme") already carried, now becoming the presentation's
# active design color scheme.
$officeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Colors($i).RGB = HexToComRgb $officeColors[$i - 1]
}
